# Fruta / hortaliza, semanal
# A new weekly record (Espinaca @ Terminal Hortofrutícola Agro Chillán) was
# added; it belongs above the existing rows 9-11, so those rows shift down
# to rows 10-12 and a brand-new row 9 is inserted with the latest data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 9 - this pushes the old rows
# 9, 10, 11 down to 10, 11, 12 respectively, preserving their values/styles.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly observation.
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value = "Ñuble"
$ws.Range("D9").Value = 44806
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = 100112012
$ws.Range("G9").Value = "Espinaca"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 120
$ws.Range("K9").Value = 7000
$ws.Range("L9").Value = 7500
$ws.Range("M9").Value = 7250
$ws.Range("N9").Value = "$/cuna 10 kilos"
$ws.Range("O9").Value = "Provincia de Diguillín"
$ws.Range("P9").Value = 725
$ws.Range("Q9").Value = 10
$ws.Range("R9").Value = "Hortaliza"

# Match the date cell's number format to the sibling rows (the insert
# should already have copied it, but make it explicit/robust).
$ws.Range("D9").NumberFormat = $ws.Range("D10").NumberFormat
